$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the page margins saved by Excel's own default template (0.7"/0.75"
# sides, 0.3" header/footer) instead of the legacy 0.75"/1"/0.5" margins
# the sheet was created with.
$ps = $ws.PageSetup
$ps.LeftMargin = 0.7 * 72
$ps.RightMargin = 0.7 * 72
$ps.TopMargin = 0.75 * 72
$ps.BottomMargin = 0.75 * 72
$ps.HeaderMargin = 0.3 * 72
$ps.FooterMargin = 0.3 * 72

# Cells whose text looks like a plain number need the column forced to
# Text first, otherwise Excel auto-coerces the input into a numeric cell
# instead of keeping it as a shared string.
$numericLookingRange = $ws.Range("C2:F2")
$numericLookingRange.NumberFormat = "@"
$ws.Range("H2").NumberFormat = "@"

$ws.Range("A2").Value = "GenCor2019"
$ws.Range("B2").Value = "AR-X"
$ws.Range("C2").Value = "32"
$ws.Range("D2").Value = "1"
$ws.Range("E2").Value = "25"
$ws.Range("F2").Value = "180"
$ws.Range("H2").Value = "1207"

# Restore the default (unstyled) cell style now that the text values are
# committed, so the cells don't carry a leftover explicit style index.
$numericLookingRange.Style = "Normal"
$ws.Range("H2").Style = "Normal"

$ws.Range("G2").ClearContents()
